# Auto-generated script applying cached-value updates (market price refresh)
# from the scheduled runner's source diff, across all 8 profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 518.25
$ws.Range("I9").Value = 491.16666
$ws.Range("K9").Value = 491.16666
$ws.Range("M9").Value = -322.16666
$ws.Range("H62").Value = 3117.6428
$ws.Range("I62").Value = 2572.111
$ws.Range("J62").Value = 4099.6
$ws.Range("K62").Value = 2572.111
$ws.Range("L62").Value = 4099.6
$ws.Range("M62").Value = -1948.111
$ws.Range("N62").Value = -5347.6
$ws.Range("H65").Value = 3117.6428
$ws.Range("I65").Value = 2572.111
$ws.Range("J65").Value = 4099.6
$ws.Range("K65").Value = 12860.555
$ws.Range("L65").Value = 20498
$ws.Range("M65").Value = -9740.555
$ws.Range("N65").Value = -26738
$ws.Range("H100").Value = 971.7778
$ws.Range("J100").Value = 633
$ws.Range("L100").Value = 633
$ws.Range("N100").Value = -1715
$ws.Range("H132").Value = 12810.139
$ws.Range("I132").Value = 2282.6177
$ws.Range("K132").Value = 6847.853099999999
$ws.Range("M132").Value = -4317.853099999999
$ws.Range("H137").Value = 290481.16
$ws.Range("I137").Value = 1001783.3
$ws.Range("J137").Value = 5960.28
$ws.Range("K137").Value = 3005349.9
$ws.Range("L137").Value = 17880.84
$ws.Range("M137").Value = -3002799.9
$ws.Range("N137").Value = -22980.84
$ws.Range("H138").Value = 3347.1282
$ws.Range("I138").Value = 2233
$ws.Range("J138").Value = 3474.457
$ws.Range("K138").Value = 6699
$ws.Range("L138").Value = 10423.371
$ws.Range("M138").Value = -1559
$ws.Range("N138").Value = -20703.371

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 349.375
$ws.Range("I5").Value = 349.375
$ws.Range("K5").Value = 349.375
$ws.Range("M5").Value = -237.375
$ws.Range("H32").Value = 4075.484
$ws.Range("I32").Value = 2021.4706
$ws.Range("K32").Value = 2021.4706
$ws.Range("M32").Value = -1734.4706
$ws.Range("H74").Value = 22729544
$ws.Range("I74").Value = 35715784
$ws.Range("K74").Value = 35715784
$ws.Range("M74").Value = -35714910
$ws.Range("H77").Value = 22729544
$ws.Range("I77").Value = 35715784
$ws.Range("K77").Value = 178578920
$ws.Range("M77").Value = -178574552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 349.375
$ws.Range("I4").Value = 349.375
$ws.Range("K4").Value = 349.375
$ws.Range("M4").Value = -234.375
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H60").Value = 75570
$ws.Range("J60").Value = 84093.336
$ws.Range("L60").Value = 84093.336
$ws.Range("N60").Value = -85291.336
$ws.Range("H80").Value = 828.625
$ws.Range("I80").Value = 1123
$ws.Range("J80").Value = 730.5
$ws.Range("K80").Value = 1123
$ws.Range("L80").Value = 730.5
$ws.Range("M80").Value = -125
$ws.Range("N80").Value = -2726.5
$ws.Range("H83").Value = 828.625
$ws.Range("I83").Value = 1123
$ws.Range("J83").Value = 730.5
$ws.Range("K83").Value = 5615
$ws.Range("L83").Value = 3652.5
$ws.Range("M83").Value = -623
$ws.Range("N83").Value = -13636.5
$ws.Range("H94").Value = 806724.7
$ws.Range("I94").Value = 914020.7
$ws.Range("K94").Value = 914020.7
$ws.Range("M94").Value = -913569.7
$ws.Range("H107").Value = 2033.1666
$ws.Range("I107").Value = 1999.75
$ws.Range("K107").Value = 1999.75
$ws.Range("M107").Value = -79.75
$ws.Range("H118").Value = 54995
$ws.Range("I118").Value = 54995
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 54995
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -53338
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").ClearContents()
$ws.Range("H123").Value = 57900
$ws.Range("J123").Value = 57900
$ws.Range("L123").Value = 57900
$ws.Range("N123").Value = -67700
$ws.Range("H133").Value = 106172
$ws.Range("J133").Value = 106172
$ws.Range("L133").Value = 106172
$ws.Range("N133").Value = -116292
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 102500
$ws.Range("J141").Value = 102500
$ws.Range("L141").Value = 102500
$ws.Range("N141").Value = -112860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21745184
$ws.Range("I31").Value = 76924610
$ws.Range("J31").Value = 7833.515
$ws.Range("K31").Value = 76924610
$ws.Range("L31").Value = 7833.515
$ws.Range("M31").Value = -76924315
$ws.Range("N31").Value = -8423.514999999999
$ws.Range("H34").Value = 21745184
$ws.Range("I34").Value = 76924610
$ws.Range("J34").Value = 7833.515
$ws.Range("K34").Value = 76924610
$ws.Range("L34").Value = 7833.515
$ws.Range("M34").Value = -76924408
$ws.Range("N34").Value = -8237.514999999999
$ws.Range("H107").Value = 998.5
$ws.Range("J107").Value = 997.5
$ws.Range("L107").Value = 997.5
$ws.Range("N107").Value = -4837.5
$ws.Range("H132").Value = 222289730
$ws.Range("I132").Value = 444511100
$ws.Range("J132").Value = 68333.336
$ws.Range("K132").Value = 1333533300
$ws.Range("L132").Value = 205000.008
$ws.Range("M132").Value = -1333530770
$ws.Range("N132").Value = -210060.008
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12465
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 818.2727
$ws.Range("I7").Value = 78.59999999999999
$ws.Range("J7").Value = 1434.6666
$ws.Range("K7").Value = 235.8
$ws.Range("L7").Value = 4303.9998
$ws.Range("M7").Value = -123.8
$ws.Range("N7").Value = -4527.9998
$ws.Range("H50").Value = 642.7143
$ws.Range("I50").Value = 479.8
$ws.Range("J50").Value = 1050
$ws.Range("K50").Value = 1439.4
$ws.Range("L50").Value = 3150
$ws.Range("M50").Value = -958.4000000000001
$ws.Range("N50").Value = -4112
$ws.Range("H53").Value = 642.7143
$ws.Range("I53").Value = 479.8
$ws.Range("J53").Value = 1050
$ws.Range("K53").Value = 1439.4
$ws.Range("L53").Value = 3150
$ws.Range("M53").Value = -958.4000000000001
$ws.Range("N53").Value = -4112
$ws.Range("H68").Value = 1724.1666
$ws.Range("J68").Value = 1724.1666
$ws.Range("L68").Value = 5172.4998
$ws.Range("N68").Value = -6794.4998
$ws.Range("H71").Value = 1724.1666
$ws.Range("J71").Value = 1724.1666
$ws.Range("L71").Value = 15517.4994
$ws.Range("N71").Value = -23629.4994
$ws.Range("H92").Value = 1023.375
$ws.Range("I92").Value = 825
$ws.Range("J92").Value = 1089.5
$ws.Range("K92").Value = 2475
$ws.Range("L92").Value = 3268.5
$ws.Range("M92").Value = -1227
$ws.Range("N92").Value = -5764.5
$ws.Range("H140").Value = 3600.879
$ws.Range("I140").Value = 2745.3635
$ws.Range("K140").Value = 8236.0905
$ws.Range("M140").Value = -3056.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 44251
$ws.Range("J93").Value = 44251
$ws.Range("L93").Value = 44251
$ws.Range("N93").Value = -47995
$ws.Range("H126").Value = 6021.7856
$ws.Range("I126").Value = 3215.6667
$ws.Range("K126").Value = 9647.000100000001
$ws.Range("M126").Value = -7177.000100000001
$ws.Range("H132").Value = 65928.94
$ws.Range("I132").Value = 108164.055
$ws.Range("J132").Value = 4200.6924
$ws.Range("K132").Value = 324492.165
$ws.Range("L132").Value = 12602.0772
$ws.Range("M132").Value = -321962.165
$ws.Range("N132").Value = -17662.0772
$ws.Range("H140").Value = 79700
$ws.Range("J140").Value = 79700
$ws.Range("L140").Value = 79700
$ws.Range("N140").Value = -90060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1098.4348
$ws.Range("I22").Value = 991.8
$ws.Range("J22").Value = 1298.375
$ws.Range("K22").Value = 991.8
$ws.Range("L22").Value = 1298.375
$ws.Range("M22").Value = -696.8
$ws.Range("N22").Value = -1888.375
$ws.Range("H27").Value = 1098.4348
$ws.Range("I27").Value = 991.8
$ws.Range("J27").Value = 1298.375
$ws.Range("K27").Value = 991.8
$ws.Range("L27").Value = 1298.375
$ws.Range("M27").Value = -884.8
$ws.Range("N27").Value = -1512.375
$ws.Range("H46").Value = 7388.1113
$ws.Range("J46").Value = 7388.1113
$ws.Range("L46").Value = 7388.1113
$ws.Range("N46").Value = -7764.1113
$ws.Range("H132").Value = 3562.3
$ws.Range("I132").Value = 3600.7083
$ws.Range("J132").Value = 3408.6667
$ws.Range("K132").Value = 10802.1249
$ws.Range("L132").Value = 10226.0001
$ws.Range("M132").Value = -8272.124899999999
$ws.Range("N132").Value = -15286.0001
$ws.Range("H134").Value = 87928.25
$ws.Range("J134").Value = 87928.25
$ws.Range("L134").Value = 87928.25
$ws.Range("N134").Value = -98068.25
$ws.Range("H138").Value = 93571
$ws.Range("J138").Value = 93571
$ws.Range("L138").Value = 93571
$ws.Range("N138").Value = -103851

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26462180
$ws.Range("I132").Value = 6946474.5
$ws.Range("J132").Value = 38471850
$ws.Range("K132").Value = 20839423.5
$ws.Range("L132").Value = 115415550
$ws.Range("M132").Value = -20836893.5
$ws.Range("N132").Value = -115420610
